$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 439, shifting existing rows 439-535 down to 441-537.
$ws.Rows("439:440").Insert()

# New row 439: Papa, Asterix, "1a (cosecha)", 2022-07-12
$ws.Range("A439").Value = 5
$ws.Range("B439").Value = "Macroferia Regional de Talca"
$ws.Range("C439").Value = "Maule"
$ws.Range("D439").Value = 44754
$ws.Range("D439").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E439").Value = 7
$ws.Range("F439").Value = 100114001
$ws.Range("G439").Value = "Papa"
$ws.Range("H439").Value = "Asterix"
$ws.Range("I439").Value = "1a (cosecha)"
$ws.Range("J439").Value = 1200
$ws.Range("K439").Value = 9000
$ws.Range("L439").Value = 9000
$ws.Range("M439").Value = 9000
$ws.Range("N439").Value = "$/saco 25 kilos"
$ws.Range("O439").Value = "Región del Maule"
$ws.Range("P439").Value = 360
$ws.Range("Q439").Value = 25
$ws.Range("R439").Value = "Hortaliza"

# New row 440: Papa, Rosara, "1a (cosecha)", 2022-07-12
$ws.Range("A440").Value = 5
$ws.Range("B440").Value = "Macroferia Regional de Talca"
$ws.Range("C440").Value = "Maule"
$ws.Range("D440").Value = 44754
$ws.Range("D440").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E440").Value = 7
$ws.Range("F440").Value = 100114001
$ws.Range("G440").Value = "Papa"
$ws.Range("H440").Value = "Rosara"
$ws.Range("I440").Value = "1a (cosecha)"
$ws.Range("J440").Value = 1600
$ws.Range("K440").Value = 6000
$ws.Range("L440").Value = 6000
$ws.Range("M440").Value = 6000
$ws.Range("N440").Value = "$/saco 25 kilos"
$ws.Range("O440").Value = "Región de La Araucanía"
$ws.Range("P440").Value = 240
$ws.Range("Q440").Value = 25
$ws.Range("R440").Value = "Hortaliza"
